$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Mobility hours value (D3) from 285 to 672
$ws.Range("D3").Value = 672

# Update the active selection to G2 (as recorded in the saved workbook view)
$ws.Range("G2").Select()
